$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 346, shifting existing rows 346:433 down to 347:434
$ws.Rows("346:346").Insert()

# Populate the newly inserted row 346 with its data.
# Columns A, B, C, E, F, G, H, I, Q, R keep the same values as the template row
# (they match the surrounding rows' pattern - copied from the row that is now 347,
# i.e. same Mercado/Region/Category/Variety/Quality/Unit-size/Classification block).
$ws.Range("A346").Value = 9
$ws.Range("B346").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C346").Value = "Metropolitana"
$ws.Range("D346").Value = 45135
$ws.Range("E346").Value = 13
$ws.Range("F346").Value = 100112021
$ws.Range("G346").Value = "Ají"
$ws.Range("H346").Value = "Americana (o)"
$ws.Range("I346").Value = "Primera"
$ws.Range("J346").Value = 52
$ws.Range("K346").Value = 34000
$ws.Range("L346").Value = 35000
$ws.Range("M346").Value = 34500
$ws.Range("N346").Value = "`$/saco 25 kilos"
$ws.Range("O346").Value = "Provincia de Limarí"
$ws.Range("P346").Value = 1380
$ws.Range("Q346").Value = 25
$ws.Range("R346").Value = "Hortaliza"

# Match the date-number style used by the other rows' Fecha (column D) cells.
$ws.Range("D346").NumberFormat = $ws.Range("D347").NumberFormat
